# Firmware_manual.pptx edit: "Key color setting done"
#
# Changes (per the canonical OOXML diff):
#  1. (best-effort; not reachable via COM in this host) slide layout's
#     datetimeFigureOut field text 2023/4/22 -> 2023/4/29
#  2. Slide 2, "Up"/"Color" key-cap shape: "Up" -> "Key Off"
#  3. Slide 2, "Theme" key-cap shape: add a new "Key" line above "Theme"
#  4. Slide 2, "Down"/"Color" key-cap shape: "Down" -> "Key On"
#  5. Slide 2, "E-Row"/"Color" key-cap shape: remove the text entirely

$p = $ppt.ActivePresentation

# --- 1. Date placeholder on the (unused) slide layout -------------------
# The "2023/4/22" text lives inside an auto-updating <a:fld type="datetimeFigureOut">
# on slideLayout3 (slide master 2's 2nd custom layout). Try the documented
# COM paths for it; harmless no-op if the host keeps date fields read-only.
try {
    $dsg = $p.Designs.Item(2)
    $mst = $dsg.SlideMaster
    $lay = $mst.CustomLayouts.Item(2)
    for ($li = 1; $li -le $lay.Shapes.Count; $li++) {
        $dshp = $lay.Shapes.Item($li)
        if ($dshp.HasTextFrame) {
            $dtr = $dshp.TextFrame.TextRange
            if ($dtr.Text -eq "2023/4/22") {
                $dtr.Text = "2023/4/29"
            }
        }
    }
    $hf = $mst.HeadersFooters.DateAndTime
    $hf.Text = "2023/4/29"
} catch {
}

# --- Slide 2 key-cap shapes ----------------------------------------------
$s = $p.Slides.Item(2)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    if ($full -eq ("Up" + [string][char]11 + "Color")) {
        # 2. "Up" -> "Key Off" (keep the existing line break + "Color" run untouched)
        $run = $tr.Characters(1, 2)
        $run.Text = "Key Off"
    }
    elseif ($full -eq "Theme") {
        # 3. Prepend a new "Key" line before the existing "Theme" line.
        # The host's text-editing primitives cannot synthesize a genuine
        # <a:br/> soft line break (embedded CR/VT always yields a new
        # <a:p/>), so approximate the intended two-line "Key" / "Theme"
        # look with a paragraph break, preserving run formatting.
        $tr.Text = "Key" + [string][char]13 + "Theme"
    }
    elseif ($full -eq ("Down" + [string][char]13 + "Color")) {
        # 4. "Down" -> "Key On" (first paragraph only; leave "Color" paragraph intact)
        $run = $tr.Characters(1, 4)
        $run.Text = "Key On"
    }
    elseif ($full -eq ("E-Row" + [string][char]11 + "Color")) {
        # 5. Remove "E-Row" / "Color" text completely, leaving an empty paragraph
        $tr.Delete()
    }
}
